$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.324.28"
$ws.Range("E2").Value = "  +1.12%  "

# Row 3
$ws.Range("D3").Value = "3.355.82"
$ws.Range("E3").Value = "  +0.95%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'584.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "

# Row 6
$ws.Range("D6").Value = "'177.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.29%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").Value = "'0.590"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.22%  "

# Row 9
$ws.Range("E9").Value = "  +3.74%  "

# Row 10
$ws.Range("E10").Value = "  +0.86%  "

# Row 11
$ws.Range("D11").Value = "'47.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.10%  "

# Row 12
$ws.Range("E12").Value = "  +1.95%  "

# Row 13
$ws.Range("D13").Value = "'690.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.64%  "

# Row 14
$ws.Range("D14").Value = "3.898.85"
$ws.Range("E14").Value = "  +0.79%  "

# Row 15
$ws.Range("D15").Value = "'8.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.83%  "

# Row 16
$ws.Range("D16").Value = "68.317.89"
$ws.Range("E16").Value = "  +1.07%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.383.81"
$ws.Range("E17").Value = "  +1.88%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.120"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.38%  "

# Row 19
$ws.Range("D19").Value = "'17.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "

# Row 20
$ws.Range("D20").Value = "'11.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.32%  "

# Row 21
$ws.Range("E21").Value = "  +0.84%  "

# Row 22
$ws.Range("D22").Value = "'5.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.86%  "

# Row 23
$ws.Range("D23").Value = "'16.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "

# Row 24
$ws.Range("D24").Value = "'100.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "

# Row 25
$ws.Range("D25").Value = "'3.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.73%  "

# Row 26
$ws.Range("E26").Value = "  +1.53%  "

# Row 27
$ws.Range("D27").Value = "'9.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.61%  "

# Row 28
$ws.Range("D28").Value = "'33.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.29%  "

# Row 29
$ws.Range("D29").Value = "'8.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.47%  "

# Row 30
$ws.Range("D30").Value = "'6.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.80%  "

# Row 31
$ws.Range("D31").Value = "'11.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.48%  "

# Row 32
$ws.Range("D32").Value = "'552.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.98%  "

# Row 33
$ws.Range("E33").Value = "  +0.59%  "

# Row 34
$ws.Range("D34").Value = "'58.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.52%  "

# Row 35
$ws.Range("D35").Value = "3.721.21"
$ws.Range("E35").Value = "  +1.23%  "

# Row 37
$ws.Range("D37").Value = "'3.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.97%  "

# Row 38
$ws.Range("E38").Value = "  +7.93%  "

# Row 39
$ws.Range("D39").Value = "'34.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.59%  "

# Row 40
$ws.Range("E40").Value = "  +2.35%  "

# Row 41
$ws.Range("D41").Value = "'2.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.25%  "

# Row 42
$ws.Range("D42").Value = "0.0ₓ0672"
$ws.Range("E42").Value = "  +1.31%  "

# Row 43
$ws.Range("D43").Value = "'0.335"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.71%  "

# Row 44
$ws.Range("D44").Value = "'3.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.26%  "

# Row 45
$ws.Range("D45").Value = "'0.0412"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.86%  "

# Row 46
$ws.Range("E46").Value = "  +2.13%  "

# Row 47
$ws.Range("E47").Value = "  +0.44%  "

# Row 48
$ws.Range("E48").Value = "  -0.05%  "

# Row 49
$ws.Range("E49").Value = "  -1.01%  "

# Row 50
$ws.Range("D50").Value = "'131.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.86%  "

# Row 51
$ws.Range("E51").Value = "  -1.42%  "
